$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a brand-new worksheet for the "2022-Q4" data, positioned right
#    after the "总计" (summary) sheet and before the existing "2021-Q2" sheet.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)
$wsNew = $wb.Worksheets.Add()
$wsNew.Move($null, $wsTotal)
$wsNew.Name = "2022-Q4"

# ---------------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a new row for 2022-Q4 above the
#    existing 2021-Q2 row (shifting it down to row 3).
# ---------------------------------------------------------------------------
# Copy formatting of the existing data row (A2:D2) down into row 3 first, so
# that the shifted 2021-Q2 row keeps its original look (notably A's style).
$wsTotal.Range("A2:D2").Copy()
$wsTotal.Range("A3:D3").PasteSpecial(-4122)

# Move the 2021-Q2 values down to row 3.
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2021-Q2"
$wsTotal.Range("C3").Value = 2
$wsTotal.Range("D3").Value = 0.03

# Write the new 2022-Q4 values into row 2.
$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 7
$wsTotal.Range("D2").Value = 0.93

# ---------------------------------------------------------------------------
# 3. Populate the new "2022-Q4" worksheet with the fund-holdings table.
# ---------------------------------------------------------------------------
$wsNew.Range("B1").Value = "基金代码"
$wsNew.Range("C1").Value = "基金名称"
$wsNew.Range("D1").Value = "基金规模"
$wsNew.Range("E1").Value = "股票总仓位"
$wsNew.Range("F1").Value = "仓位占比"
$wsNew.Range("G1").Value = "持有市值(亿元)"
$wsNew.Range("H1").Value = "仓位排名"

$rows = @(
    @{ A=0; B="007139"; C="富国民裕进取沪港深成长精选混合A"; D="13.04"; E="92.50"; F="4.85"; G="0.6324"; H=9  },
    @{ A=1; B="011556"; C="富国民裕进取沪港深成长精选混合C"; D="2.60";  E="92.50"; F="4.85"; G="0.1261"; H=9  },
    @{ A=2; B="004317"; C="前海开源沪港深裕鑫灵活配置混合C"; D="2.88";  E="90.85"; F="3.06"; G="0.0881"; H=6  },
    @{ A=3; B="004316"; C="前海开源沪港深裕鑫灵活配置混合A"; D="2.30";  E="90.85"; F="3.06"; G="0.0704"; H=6  },
    @{ A=4; B="006106"; C="景顺长城量化港股通股票";         D="0.55";  E="81.25"; F="1.33"; G="0.0073"; H=10 },
    @{ A=5; B="013989"; C="富国沪港深优质资产混合A";         D="0.14";  E="90.56"; F="3.81"; G="0.0053"; H=6  },
    @{ A=6; B="013990"; C="富国沪港深优质资产混合C";         D="0.05";  E="90.56"; F="3.81"; G="0.0019"; H=6  }
)

$r = 2
foreach ($row in $rows) {
    $wsNew.Range("A$r").Value = $row.A
    # Columns B, D, E, F, G hold numeric-looking text (fund codes with
    # leading zeros, e.g. "007139", and decimals like "13.04") in the
    # source data, so a leading apostrophe forces Excel to store them as
    # text instead of auto-converting to numbers (which would drop the
    # leading zeros / trailing zeros).
    $wsNew.Range("B$r").Value = "'" + $row.B
    $wsNew.Range("C$r").Value = $row.C
    $wsNew.Range("D$r").Value = "'" + $row.D
    $wsNew.Range("E$r").Value = "'" + $row.E
    $wsNew.Range("F$r").Value = "'" + $row.F
    $wsNew.Range("G$r").Value = "'" + $row.G
    $wsNew.Range("H$r").Value = $row.H
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 4. Match styling: the header row and column-A cells on the new sheet use
#    the same bold/bordered style already used by the "总计" header (B1).
# ---------------------------------------------------------------------------
$wsTotal.Range("B1").Copy()
$wsNew.Range("B1:H1").PasteSpecial(-4122)

$wsTotal.Range("A2").Copy()
$wsNew.Range("A2:A8").PasteSpecial(-4122)

# The leading apostrophe used above to force text storage also marks the
# cells with a "quote prefix" style flag. Drop that formatting residue so
# the cells end up with no explicit style, matching the source data (which
# keeps B/D/E/F/G as plain, unstyled text cells).
$wsNew.Range("B2:B8").ClearFormats()
$wsNew.Range("D2:G8").ClearFormats()

Write-Host "done"
